# edit.ps1 - apply Catanzaro.xlsx update: extend each sheet with new daily rows
# (through 2021-07-31 / serial 44408), adding case counts for 2021-06-27 to 2021-07-06
# (serials 44374-44382) and the corresponding rolling 7-day average formulas.
$wb = $excel.ActiveWorkbook

# ---- Sheet 1: Nuovi casi ----
$ws = $wb.Worksheets.Item(1)

# Rows 476-484: date in col A, new count in col C, rolling 7-day average in col D
$a = $ws.Cells.Item(476, 1)
$a.Value = 44374
$a.NumberFormat = "dd/mm/yyyy"
$c = $ws.Cells.Item(476, 3)
$c.Value = 0
$ws.Cells.Item(476, 4).Formula = "=AVERAGE(C470:C476)"
$a = $ws.Cells.Item(477, 1)
$a.Value = 44375
$a.NumberFormat = "dd/mm/yyyy"
$c = $ws.Cells.Item(477, 3)
$c.Value = 0
$ws.Cells.Item(477, 4).Formula = "=AVERAGE(C471:C477)"
$a = $ws.Cells.Item(478, 1)
$a.Value = 44376
$a.NumberFormat = "dd/mm/yyyy"
$c = $ws.Cells.Item(478, 3)
$c.Value = 2
$ws.Cells.Item(478, 4).Formula = "=AVERAGE(C472:C478)"
$a = $ws.Cells.Item(479, 1)
$a.Value = 44377
$a.NumberFormat = "dd/mm/yyyy"
$c = $ws.Cells.Item(479, 3)
$c.Value = 4
$ws.Cells.Item(479, 4).Formula = "=AVERAGE(C473:C479)"
$a = $ws.Cells.Item(480, 1)
$a.Value = 44378
$a.NumberFormat = "dd/mm/yyyy"
$c = $ws.Cells.Item(480, 3)
$c.Value = 0
$ws.Cells.Item(480, 4).Formula = "=AVERAGE(C474:C480)"
$a = $ws.Cells.Item(481, 1)
$a.Value = 44379
$a.NumberFormat = "dd/mm/yyyy"
$c = $ws.Cells.Item(481, 3)
$c.Value = 6
$ws.Cells.Item(481, 4).Formula = "=AVERAGE(C475:C481)"
$a = $ws.Cells.Item(482, 1)
$a.Value = 44380
$a.NumberFormat = "dd/mm/yyyy"
$c = $ws.Cells.Item(482, 3)
$c.Value = 1
$ws.Cells.Item(482, 4).Formula = "=AVERAGE(C476:C482)"
$a = $ws.Cells.Item(483, 1)
$a.Value = 44381
$a.NumberFormat = "dd/mm/yyyy"
$c = $ws.Cells.Item(483, 3)
$c.Value = 0
$ws.Cells.Item(483, 4).Formula = "=AVERAGE(C477:C483)"
$a = $ws.Cells.Item(484, 1)
$a.Value = 44382
$a.NumberFormat = "dd/mm/yyyy"
$c = $ws.Cells.Item(484, 3)
$c.Value = 0
$ws.Cells.Item(484, 4).Formula = "=AVERAGE(C478:C484)"

# Rows 485-510: date column only (no case counts available yet for these future days)
for ($r = 485; $r -le 510; $r++) {
    $a = $ws.Cells.Item($r, 1)
    $a.Value = 44374 + ($r - 476)
    $a.NumberFormat = "dd/mm/yyyy"
}

# ---- Sheet 2: Deceduti ----
$ws = $wb.Worksheets.Item(2)

# Rows 476-484: date in col A, new count in col C, rolling 7-day average in col D
$a = $ws.Cells.Item(476, 1)
$a.Value = 44374
$a.NumberFormat = "dd/mm/yyyy"
$c = $ws.Cells.Item(476, 3)
$c.Value = 0
$ws.Cells.Item(476, 4).Formula = "=AVERAGE(C470:C476)"
$a = $ws.Cells.Item(477, 1)
$a.Value = 44375
$a.NumberFormat = "dd/mm/yyyy"
$c = $ws.Cells.Item(477, 3)
$c.Value = 0
$ws.Cells.Item(477, 4).Formula = "=AVERAGE(C471:C477)"
$a = $ws.Cells.Item(478, 1)
$a.Value = 44376
$a.NumberFormat = "dd/mm/yyyy"
$c = $ws.Cells.Item(478, 3)
$c.Value = 0
$ws.Cells.Item(478, 4).Formula = "=AVERAGE(C472:C478)"
$a = $ws.Cells.Item(479, 1)
$a.Value = 44377
$a.NumberFormat = "dd/mm/yyyy"
$c = $ws.Cells.Item(479, 3)
$c.Value = 0
$ws.Cells.Item(479, 4).Formula = "=AVERAGE(C473:C479)"
$a = $ws.Cells.Item(480, 1)
$a.Value = 44378
$a.NumberFormat = "dd/mm/yyyy"
$c = $ws.Cells.Item(480, 3)
$c.Value = 1
$ws.Cells.Item(480, 4).Formula = "=AVERAGE(C474:C480)"
$a = $ws.Cells.Item(481, 1)
$a.Value = 44379
$a.NumberFormat = "dd/mm/yyyy"
$c = $ws.Cells.Item(481, 3)
$c.Value = 0
$ws.Cells.Item(481, 4).Formula = "=AVERAGE(C475:C481)"
$a = $ws.Cells.Item(482, 1)
$a.Value = 44380
$a.NumberFormat = "dd/mm/yyyy"
$c = $ws.Cells.Item(482, 3)
$c.Value = 0
$ws.Cells.Item(482, 4).Formula = "=AVERAGE(C476:C482)"
$a = $ws.Cells.Item(483, 1)
$a.Value = 44381
$a.NumberFormat = "dd/mm/yyyy"
$c = $ws.Cells.Item(483, 3)
$c.Value = 0
$ws.Cells.Item(483, 4).Formula = "=AVERAGE(C477:C483)"
$a = $ws.Cells.Item(484, 1)
$a.Value = 44382
$a.NumberFormat = "dd/mm/yyyy"
$c = $ws.Cells.Item(484, 3)
$c.Value = 0
$ws.Cells.Item(484, 4).Formula = "=AVERAGE(C478:C484)"

# Rows 485-510: date column only (no case counts available yet for these future days)
for ($r = 485; $r -le 510; $r++) {
    $a = $ws.Cells.Item($r, 1)
    $a.Value = 44374 + ($r - 476)
    $a.NumberFormat = "dd/mm/yyyy"
}

# ---- Sheet 3: Dimessi   Guariti ----
$ws = $wb.Worksheets.Item(3)

# Rows 476-484: date in col A, new count in col C, rolling 7-day average in col D
$a = $ws.Cells.Item(476, 1)
$a.Value = 44374
$a.NumberFormat = "dd/mm/yyyy"
$c = $ws.Cells.Item(476, 3)
$c.Value = 2
$ws.Cells.Item(476, 4).Formula = "=AVERAGE(C470:C476)"
$a = $ws.Cells.Item(477, 1)
$a.Value = 44375
$a.NumberFormat = "dd/mm/yyyy"
$c = $ws.Cells.Item(477, 3)
$c.Value = 29
$ws.Cells.Item(477, 4).Formula = "=AVERAGE(C471:C477)"
$a = $ws.Cells.Item(478, 1)
$a.Value = 44376
$a.NumberFormat = "dd/mm/yyyy"
$c = $ws.Cells.Item(478, 3)
$c.Value = 11
$ws.Cells.Item(478, 4).Formula = "=AVERAGE(C472:C478)"
$a = $ws.Cells.Item(479, 1)
$a.Value = 44377
$a.NumberFormat = "dd/mm/yyyy"
$c = $ws.Cells.Item(479, 3)
$c.Value = 37
$ws.Cells.Item(479, 4).Formula = "=AVERAGE(C473:C479)"
$a = $ws.Cells.Item(480, 1)
$a.Value = 44378
$a.NumberFormat = "dd/mm/yyyy"
$c = $ws.Cells.Item(480, 3)
$c.Value = 11
$ws.Cells.Item(480, 4).Formula = "=AVERAGE(C474:C480)"
$a = $ws.Cells.Item(481, 1)
$a.Value = 44379
$a.NumberFormat = "dd/mm/yyyy"
$c = $ws.Cells.Item(481, 3)
$c.Value = 41
$ws.Cells.Item(481, 4).Formula = "=AVERAGE(C475:C481)"
$a = $ws.Cells.Item(482, 1)
$a.Value = 44380
$a.NumberFormat = "dd/mm/yyyy"
$c = $ws.Cells.Item(482, 3)
$c.Value = 0
$ws.Cells.Item(482, 4).Formula = "=AVERAGE(C476:C482)"
$a = $ws.Cells.Item(483, 1)
$a.Value = 44381
$a.NumberFormat = "dd/mm/yyyy"
$c = $ws.Cells.Item(483, 3)
$c.Value = 0
$ws.Cells.Item(483, 4).Formula = "=AVERAGE(C477:C483)"
$a = $ws.Cells.Item(484, 1)
$a.Value = 44382
$a.NumberFormat = "dd/mm/yyyy"
$c = $ws.Cells.Item(484, 3)
$c.Value = 25
$ws.Cells.Item(484, 4).Formula = "=AVERAGE(C478:C484)"

# Rows 485-510: date column only (no case counts available yet for these future days)
for ($r = 485; $r -le 510; $r++) {
    $a = $ws.Cells.Item($r, 1)
    $a.Value = 44374 + ($r - 476)
    $a.NumberFormat = "dd/mm/yyyy"
}

# ---- Sheet 4: Ricoveri ----
$ws = $wb.Worksheets.Item(4)

# Rows 476-484: date in col A, new count in col C, rolling 7-day average in col D
$a = $ws.Cells.Item(476, 1)
$a.Value = 44374
$a.NumberFormat = "dd/mm/yyyy"
$c = $ws.Cells.Item(476, 3)
$c.Value = 14
$c.Font.Color = 0
$ws.Cells.Item(476, 4).Formula = "=AVERAGE(C470:C476)"
$a = $ws.Cells.Item(477, 1)
$a.Value = 44375
$a.NumberFormat = "dd/mm/yyyy"
$c = $ws.Cells.Item(477, 3)
$c.Value = 15
$c.Font.Color = 0
$ws.Cells.Item(477, 4).Formula = "=AVERAGE(C471:C477)"
$a = $ws.Cells.Item(478, 1)
$a.Value = 44376
$a.NumberFormat = "dd/mm/yyyy"
$c = $ws.Cells.Item(478, 3)
$c.Value = 16
$c.Font.Color = 0
$ws.Cells.Item(478, 4).Formula = "=AVERAGE(C472:C478)"
$a = $ws.Cells.Item(479, 1)
$a.Value = 44377
$a.NumberFormat = "dd/mm/yyyy"
$c = $ws.Cells.Item(479, 3)
$c.Value = 16
$c.Font.Color = 0
$ws.Cells.Item(479, 4).Formula = "=AVERAGE(C473:C479)"
$a = $ws.Cells.Item(480, 1)
$a.Value = 44378
$a.NumberFormat = "dd/mm/yyyy"
$c = $ws.Cells.Item(480, 3)
$c.Value = 14
$c.Font.Color = 0
$ws.Cells.Item(480, 4).Formula = "=AVERAGE(C474:C480)"
$a = $ws.Cells.Item(481, 1)
$a.Value = 44379
$a.NumberFormat = "dd/mm/yyyy"
$c = $ws.Cells.Item(481, 3)
$c.Value = 11
$c.Font.Color = 0
$ws.Cells.Item(481, 4).Formula = "=AVERAGE(C475:C481)"
$a = $ws.Cells.Item(482, 1)
$a.Value = 44380
$a.NumberFormat = "dd/mm/yyyy"
$c = $ws.Cells.Item(482, 3)
$c.Value = 12
$c.Font.Color = 0
$ws.Cells.Item(482, 4).Formula = "=AVERAGE(C476:C482)"
$a = $ws.Cells.Item(483, 1)
$a.Value = 44381
$a.NumberFormat = "dd/mm/yyyy"
$c = $ws.Cells.Item(483, 3)
$c.Value = 12
$c.Font.Color = 0
$ws.Cells.Item(483, 4).Formula = "=AVERAGE(C477:C483)"
$a = $ws.Cells.Item(484, 1)
$a.Value = 44382
$a.NumberFormat = "dd/mm/yyyy"
$c = $ws.Cells.Item(484, 3)
$c.Value = 10
$c.Font.Color = 0
$ws.Cells.Item(484, 4).Formula = "=AVERAGE(C478:C484)"

# Rows 485-510: date column only (no case counts available yet for these future days)
for ($r = 485; $r -le 510; $r++) {
    $a = $ws.Cells.Item($r, 1)
    $a.Value = 44374 + ($r - 476)
    $a.NumberFormat = "dd/mm/yyyy"
}

# ---- Sheet 5: Terapia ----
$ws = $wb.Worksheets.Item(5)

# Rows 476-484: date in col A, new count in col C, rolling 7-day average in col D
$a = $ws.Cells.Item(476, 1)
$a.Value = 44374
$a.NumberFormat = "dd/mm/yyyy"
$c = $ws.Cells.Item(476, 3)
$c.Value = 2
$c.Font.Color = 0
$ws.Cells.Item(476, 4).Formula = "=AVERAGE(C470:C476)"
$a = $ws.Cells.Item(477, 1)
$a.Value = 44375
$a.NumberFormat = "dd/mm/yyyy"
$c = $ws.Cells.Item(477, 3)
$c.Value = 2
$c.Font.Color = 0
$ws.Cells.Item(477, 4).Formula = "=AVERAGE(C471:C477)"
$a = $ws.Cells.Item(478, 1)
$a.Value = 44376
$a.NumberFormat = "dd/mm/yyyy"
$c = $ws.Cells.Item(478, 3)
$c.Value = 2
$c.Font.Color = 0
$ws.Cells.Item(478, 4).Formula = "=AVERAGE(C472:C478)"
$a = $ws.Cells.Item(479, 1)
$a.Value = 44377
$a.NumberFormat = "dd/mm/yyyy"
$c = $ws.Cells.Item(479, 3)
$c.Value = 1
$c.Font.Color = 0
$ws.Cells.Item(479, 4).Formula = "=AVERAGE(C473:C479)"
$a = $ws.Cells.Item(480, 1)
$a.Value = 44378
$a.NumberFormat = "dd/mm/yyyy"
$c = $ws.Cells.Item(480, 3)
$c.Value = 1
$c.Font.Color = 0
$ws.Cells.Item(480, 4).Formula = "=AVERAGE(C474:C480)"
$a = $ws.Cells.Item(481, 1)
$a.Value = 44379
$a.NumberFormat = "dd/mm/yyyy"
$c = $ws.Cells.Item(481, 3)
$c.Value = 1
$c.Font.Color = 0
$ws.Cells.Item(481, 4).Formula = "=AVERAGE(C475:C481)"
$a = $ws.Cells.Item(482, 1)
$a.Value = 44380
$a.NumberFormat = "dd/mm/yyyy"
$c = $ws.Cells.Item(482, 3)
$c.Value = 1
$c.Font.Color = 0
$ws.Cells.Item(482, 4).Formula = "=AVERAGE(C476:C482)"
$a = $ws.Cells.Item(483, 1)
$a.Value = 44381
$a.NumberFormat = "dd/mm/yyyy"
$c = $ws.Cells.Item(483, 3)
$c.Value = 1
$c.Font.Color = 0
$ws.Cells.Item(483, 4).Formula = "=AVERAGE(C477:C483)"
$a = $ws.Cells.Item(484, 1)
$a.Value = 44382
$a.NumberFormat = "dd/mm/yyyy"
$c = $ws.Cells.Item(484, 3)
$c.Value = 1
$c.Font.Color = 0
$ws.Cells.Item(484, 4).Formula = "=AVERAGE(C478:C484)"

# Rows 485-510: date column only (no case counts available yet for these future days)
for ($r = 485; $r -le 510; $r++) {
    $a = $ws.Cells.Item($r, 1)
    $a.Value = 44374 + ($r - 476)
    $a.NumberFormat = "dd/mm/yyyy"
}

# ---- View state: selection per sheet, and active tab ----

$ws = $wb.Worksheets.Item(1)
$ws.Range("A476:D484").Select()
$ws = $wb.Worksheets.Item(2)
$ws.Range("A476:D484").Select()
$ws = $wb.Worksheets.Item(3)
$ws.Range("A476:D484").Select()
$ws = $wb.Worksheets.Item(4)
$ws.Range("A476:D484").Select()
$ws = $wb.Worksheets.Item(5)
$ws.Range("A476:D484").Select()

# "Dimessi   Guariti" (sheet 3) is the active tab after the edit
$wb.Worksheets.Item(3).Activate()

